# "more gameflow, start -> level 2"
# Insert a new "Level 2" block (igneous/sedimentary/metamorphic key/value
# pairs) right before the existing igneous-rock entries, shifting all the
# rows that followed down by three.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 5-16 down to 8-19, opening up three blank rows at 5:7.
$ws.Rows("5:7").Insert() | Out-Null

# Fill column A (Key) first, then column B (Value) - matches authoring order.
$ws.Range("A5").Value = "igneous"
$ws.Range("A6").Value = "sedimentary"
$ws.Range("A7").Value = "metamorphic"

$ws.Range("B5").Value = "Igneous"
$ws.Range("B6").Value = "Sedimentary"
$ws.Range("B7").Value = "Metamorphic"

# Move the active selection to B7, matching the edited workbook's view state.
$ws.Range("B7").Select() | Out-Null
